# Generate Report for Handoff
# Update status strings ("In Translation" -> "Ready for handoff") and
# refresh the related timestamp strings, then widen the "Status" /
# per-language status columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus       = "Ready for handoff"
$newHoXliffDate  = "2016-09-05 18:44:36"
$newHandoffDate  = "2016-09-05 18:44:32"

# --- Overview sheet: status columns (E, F) and HO Xliff generate date (G) ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newHoXliffDate

# --- zh-cn sheet: Status (C) and Latest Handoff Datetime (H) ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $newHandoffDate

# --- de-de sheet: Status (C) and Latest Handoff Datetime (H) ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $newHoXliffDate

# --- Widen columns to accommodate the new "Ready for handoff" text ---
# Target stored column width is 17.2159881591797 "Excel character" units.
# The ColumnWidth setter here only lands on multiples of 1/6, so feed it a
# value whose quantized result (17.166666666666668) is the closest
# reachable approximation of the target width.
$targetColumnWidth = 16.33
$wsOverview.Range("E:F").ColumnWidth = $targetColumnWidth
$wsZhCn.Range("C:C").ColumnWidth = $targetColumnWidth
$wsDeDe.Range("C:C").ColumnWidth = $targetColumnWidth
